$d = $word.ActiveDocument

# 1. Fix the typo'd id: p0061v_1 -> p061v_1
$d.Content.Find.Execute("p0061v_1", $true, $false, $false, $false, $false, $true, 1, $false, "p061v_1", 2) | Out-Null

# 2. Collapse the split <id>...</id> runs for p061v_2, p061v_3 and p061v_4 into
#    a single run each (keeping the Courier New formatting of the <id> run).
$d.Content.Find.Execute("<id>p061v_2</id>", $true, $false, $false, $false, $false, $true, 1, $false, "<id>p061v_2</id>", 2) | Out-Null
$d.Content.Find.Execute("<id>p061v_3</id>", $true, $false, $false, $false, $false, $true, 1, $false, "<id>p061v_3</id>", 2) | Out-Null
$d.Content.Find.Execute("<id>p061v_4</id>", $true, $false, $false, $false, $false, $true, 1, $false, "<id>p061v_4</id>", 2) | Out-Null
